$d = $word.ActiveDocument

# Avoid Word's AutoFormat turning our straight quote into a curly one
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# 1) Replace "... транспортного средства N  ${document_number} " with
#    "... транспортного средства от" (text-level result of the diff).
$find1 = $d.Content
$find1.Find.Execute(
    "транспортного средства N  `$\{document_number\} ",
    $true, $false, $true, $false, $false,
    $true, 1, $false, "транспортного средства от",
    2
) | Out-Null

# 2) The "_GoBack" bookmark needs to end up collapsed right after the newly
#    inserted "от" (i.e. immediately before the following quotation mark),
#    matching where the block of runs now ends up relative to the bookmark.
$locate = $d.Content
$locate.Find.Execute('от"', $true, $false, $false) | Out-Null
$newBookmarkPos = $locate.End - 1

$target = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
